# Add 2022-Q4 data:
#  - insert a "2022-Q4" row at the top of the "总计" (totals) sheet's data
#    (pushing 2022-Q3 / 2022-Q1 / 2021-Q4 / 2021-Q2 down by one row)
#  - add a brand-new "2022-Q4" worksheet (positioned right after "总计",
#    before "2022-Q3") with the per-fund holdings detail

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force a value to be written as text even when it "looks" numeric
    # (fund codes like 005360 must keep their leading zeros; the ratio
    # columns are stored as text too), then strip the transient
    # quote-prefix formatting ClearFormats() would otherwise leave
    # behind. Only call this on cells that carry NO intentional style
    # of their own (ClearFormats wipes styling wholesale).
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift existing rows down and insert the 2022-Q4 row
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Copy the formatted row 5 down to row 6 first so the new last row keeps
# the same per-cell styling (bordered/bold index cell in column A) as
# every other data row, then overwrite every row's contents bottom-up.
$totals.Range("A5").Copy($totals.Range("A6"))

$totals.Range("A6").Value = 4
$totals.Range("B6").Value = "2021-Q2"
$totals.Range("C6").Value = 2
$totals.Range("D6").Value = 0

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2021-Q4"
$totals.Range("C5").Value = 5
$totals.Range("D5").Value = 0.22

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 17
$totals.Range("D3").Value = 1.47

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 6
$totals.Range("D2").Value = 0.3

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q4" worksheet, placed right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q4.Name = "2022-Q4"

# Borrow the exact column widths / header & index-column styling from an
# existing quarter sheet (same layout in every quarter tab) and then
# overwrite every cell's contents with the 2022-Q4 figures.
$template = $wb.Worksheets.Item("2022-Q3")
$template.Range("A1:H7").Copy($q4.Range("A1"))

# Header labels are plain (non-numeric-looking) text, so a normal
# assignment keeps them as text without disturbing the bold/border
# header style (s="2") the template copy already applied.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows: column A is the plain numeric index (keeps its s="2" style
# from the template copy); B/D/E/F/G are numeric-looking text that must
# keep leading zeros / exact decimal text; C is free-form fund-name text
# (never numeric-looking); H is a genuine number.
$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "506007"
$q4.Range("C2").Value = "广发科创板两年定开混合"
Set-TextValue $q4.Range("D2") "5.31"
Set-TextValue $q4.Range("E2") "88.81"
Set-TextValue $q4.Range("F2") "4.87"
Set-TextValue $q4.Range("G2") "0.2586"
$q4.Range("H2").Value = 7

$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "005360"
$q4.Range("C3").Value = "汇安资产轮动灵活配置混合A"
Set-TextValue $q4.Range("D3") "0.26"
Set-TextValue $q4.Range("E3") "94.03"
Set-TextValue $q4.Range("F3") "6.71"
Set-TextValue $q4.Range("G3") "0.0174"
$q4.Range("H3").Value = 7

$q4.Range("A4").Value = 2
Set-TextValue $q4.Range("B4") "620001"
$q4.Range("C4").Value = "金元顺安宝石动力混合"
Set-TextValue $q4.Range("D4") "0.46"
Set-TextValue $q4.Range("E4") "56.89"
Set-TextValue $q4.Range("F4") "3.73"
Set-TextValue $q4.Range("G4") "0.0172"
$q4.Range("H4").Value = 9

$q4.Range("A5").Value = 3
Set-TextValue $q4.Range("B5") "006231"
$q4.Range("C5").Value = "国融融君灵活配置混合A"
Set-TextValue $q4.Range("D5") "0.10"
Set-TextValue $q4.Range("E5") "59.86"
Set-TextValue $q4.Range("F5") "2.16"
Set-TextValue $q4.Range("G5") "0.0022"
$q4.Range("H5").Value = 6

$q4.Range("A6").Value = 4
Set-TextValue $q4.Range("B6") "017213"
$q4.Range("C6").Value = "汇安资产轮动灵活配置混合C"
Set-TextValue $q4.Range("D6") "0.01"
Set-TextValue $q4.Range("E6") "94.03"
Set-TextValue $q4.Range("F6") "6.71"
Set-TextValue $q4.Range("G6") "0.0007"
$q4.Range("H6").Value = 7

$q4.Range("A7").Value = 5
Set-TextValue $q4.Range("B7") "006232"
$q4.Range("C7").Value = "国融融君灵活配置混合C"
Set-TextValue $q4.Range("D7") "0.01"
Set-TextValue $q4.Range("E7") "59.86"
Set-TextValue $q4.Range("F7") "2.16"
Set-TextValue $q4.Range("G7") "0.0002"
$q4.Range("H7").Value = 6

Write-Host "2022-Q4 sheet + totals row added"
